{"js": "// The paragraph containing the M2Doc field text ends with the run\n// \"...)}\" \u2014 the closing parenthesis of the asTable(...) call immediately\n// followed by the field's closing brace \"}\". The parser migration needs\n// these two characters to live in separate runs (same formatting), i.e.\n// \"...)\" + \"}\" instead of one combined \"...)}\" run.\n//\n// There is exactly one \")\" character in the whole document, immediately\n// followed by the field-closing \"}\", so searching for the lone \")\" gives\n// us precisely the boundary we need to split on.\nconst results = context.document.body.search(\")\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \")\" to split the run on.');\n}\n\nconst closingParenRange = results.items[0];\n\n// Toggling a character-formatting property on just the \")\" range and then\n// reverting it forces the run that used to hold \")}\" to split into two\n// runs at that boundary \u2014 one for \")\" and one for \"}\" \u2014 without leaving\n// any residual direct formatting behind, since the value is restored.\nclosingParenRange.font.bold = true;\nawait context.sync();\n\nclosingParenRange.font.bold = false;\nawait context.sync();\n", "ps1": "# The paragraph containing the M2Doc field text ends with a single run\n# \"...)}\" \u2014 the closing parenthesis of the asTable(...) call immediately\n# followed by the field's closing brace \"}\". The parser migration needs\n# these two characters to live in separate runs (with identical\n# formatting), i.e. \"...)\" + \"}\" instead of one combined \"...)}\" run.\n\n$d = $word.ActiveDocument\n\n# There is exactly one \")\" character in the whole document, immediately\n# followed by the field-closing \"}\", so searching for the lone \")\" gives\n# us precisely the boundary we need to split on.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\")\")\n\nif ($rng.Find.Found) {\n    # Toggling a character-formatting property on just the \")\" character\n    # and then reverting it forces the run that used to hold \")}\" to split\n    # into two runs at that boundary \u2014 one for \")\" and one for \"}\" \u2014 with\n    # no residual direct formatting left behind, since the value is\n    # restored right after.\n    $rng.Bold = 1\n    $rng.Bold = 0\n}\n"}
